# Edit script: "practiced basic of flutter"
# 1. Append new content (page break, "Text Field" section, "State full Widgets" section)
#    at the end of the document.
# 2. Update the ObjectID of the 5 embedded OLE objects.
# 3. Remove the old "_GoBack" bookmark (the new content carries its own, which the
#    engine will renumber automatically, matching the committed document).

$d = $word.ActiveDocument

function Wrap-Body {
    param($innerXml)
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Step 1: append the new paragraphs after the current last paragraph.
# (InsertXML on the very last paragraph of the body leaves a spare empty
# paragraph behind because the final mark is shared with sectPr, so we strip
# that extra paragraph mark away again afterwards.)
# ---------------------------------------------------------------------------

$newContent = @'
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Text Field</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>The only</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> state full widget in </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>flutter</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">State </w:t></w:r><w:bookmarkStart w:id="500" w:name="_GoBack"/><w:bookmarkEnd w:id="500"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>full Widgets</w:t></w:r></w:p>
'@

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.InsertXML((Wrap-Body $newContent))

# Trim the spare trailing empty paragraph mark that appears after replacing
# the body's final paragraph.
$n = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($n - 1)
$trim = $d.Range($secondToLast.Range.End - 1, $d.Content.End)
$trim.Delete()

# ---------------------------------------------------------------------------
# Step 2: fix up the ObjectID of each embedded OLE object (paragraphs 3, 6,
# 12, 17 and 19 are untouched by the insertion above, since that happened
# after paragraph 19).
# ---------------------------------------------------------------------------

function Set-ParagraphXml {
    param($paraIndex, $newXml)
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.InsertXML((Wrap-Body $newXml))
}

$para3 = @'
<w:p w:rsidR="00162D26" w:rsidRDefault="00162D26" w:rsidP="00162D26"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:object w:dxaOrig="9026" w:dyaOrig="10830"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:451.5pt;height:541.5pt" o:ole=""><v:imagedata r:id="rId4" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Word.OpenDocumentText.12" ShapeID="_x0000_i1025" DrawAspect="Content" ObjectID="_1732088574" r:id="rId5"/></w:object></w:r></w:p>
'@
Set-ParagraphXml 3 $para3

$para6 = @'
<w:p w:rsidR="00162D26" w:rsidRDefault="004527AD" w:rsidP="00162D26"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:object w:dxaOrig="9026" w:dyaOrig="4560"><v:shape id="_x0000_i1026" type="#_x0000_t75" style="width:451.5pt;height:228pt" o:ole=""><v:imagedata r:id="rId6" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Word.OpenDocumentText.12" ShapeID="_x0000_i1026" DrawAspect="Content" ObjectID="_1732088575" r:id="rId7"/></w:object></w:r></w:p>
'@
Set-ParagraphXml 6 $para6

$para12 = @'
<w:p w:rsidR="00864CC0" w:rsidRDefault="00A7001B" w:rsidP="00864CC0"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:object w:dxaOrig="9026" w:dyaOrig="4275"><v:shape id="_x0000_i1027" type="#_x0000_t75" style="width:451.5pt;height:214pt" o:ole=""><v:imagedata r:id="rId8" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Word.OpenDocumentText.12" ShapeID="_x0000_i1027" DrawAspect="Content" ObjectID="_1732088576" r:id="rId9"/></w:object></w:r></w:p>
'@
Set-ParagraphXml 12 $para12

$para17 = @'
<w:p w:rsidR="008E4692" w:rsidRDefault="00D73397" w:rsidP="00ED0861"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:object w:dxaOrig="9026" w:dyaOrig="4560"><v:shape id="_x0000_i1028" type="#_x0000_t75" style="width:451.5pt;height:228pt" o:ole=""><v:imagedata r:id="rId10" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Word.OpenDocumentText.12" ShapeID="_x0000_i1028" DrawAspect="Content" ObjectID="_1732088577" r:id="rId11"/></w:object></w:r></w:p>
'@
Set-ParagraphXml 17 $para17

$para19 = @'
<w:p w:rsidR="00D73397" w:rsidRPr="00ED0861" w:rsidRDefault="00722188" w:rsidP="00ED0861"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:object w:dxaOrig="9026" w:dyaOrig="1020"><v:shape id="_x0000_i1029" type="#_x0000_t75" style="width:451.5pt;height:51pt" o:ole=""><v:imagedata r:id="rId12" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Word.OpenDocumentText.12" ShapeID="_x0000_i1029" DrawAspect="Content" ObjectID="_1732088578" r:id="rId13"/></w:object></w:r></w:p>
'@
Set-ParagraphXml 19 $para19

# ---------------------------------------------------------------------------
# Step 3: remove the original "_GoBack" bookmark. The new "State "/"full
# Widgets" paragraph inserted above carries the new "_GoBack" bookmark, and
# the engine renumbers bookmark ids automatically once the stale one is gone.
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
